$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "vendedor" (salesperson) is now an optional field, so this new sale
# is recorded without one (and without a "cliente" either).
$ws.Cells.Item(6, 1).Value = "22/02/2023"
$ws.Cells.Item(6, 2).Value = ""
$ws.Cells.Item(6, 3).Value = ""
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = -108.01
